# Update the lab 9 title on slide 1: split the subtitle's first line into
# "The Central Limit Theorem " + "and Sample Means" (was a single run
# "Skewness, Normality, and Sample Means").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)               # "Subtitle 2" placeholder
$tr2 = $shape.TextFrame2.TextRange
$para1 = $tr2.Paragraphs(1)              # first paragraph only - leave the rest of the subtitle untouched

$newFirst  = "The Central Limit Theorem "
$newSecond = "and Sample Means"

# Rewrite the paragraph text in place, then carve it into two runs by
# re-setting the text of each character sub-range (this preserves the
# existing run formatting/rPr while forcing a run split at the boundary).
$para1.Text = $newFirst + $newSecond

$firstRange = $para1.Characters(1, $newFirst.Length)
$firstRange.Text = $newFirst

$secondRange = $para1.Characters($newFirst.Length + 1, $newSecond.Length)
$secondRange.Text = $newSecond
